$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct an existing number (Recover count for row 25) ---
$ws.Range("H25").Value = 4

# --- Append new daily records (rows 29-33), 11-15 April 2020 ---

# Row 29 : 11-Apr-2020
$ws.Range("A29").Value = 43932
$ws.Range("A29").NumberFormat = "d-mmm-yy"
$ws.Range("A29").HorizontalAlignment = -4108
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 491
$ws.Range("E29").Value = "Nairobi,Mombasa"
$ws.Range("G29").Value = "Community(2)"
$ws.Range("D29").Value = "None"

# --- Correct an existing entry (Travelled From for row 7) ---
$ws.Range("D7").Value = "Spain(2) , Dubai"

$ws.Range("F29").Value = 191
$ws.Range("H29").Value = 2
$ws.Range("I29").Value = 0
$ws.Range("L29").Value = "02 - 32."
$ws.Range("L29").NumberFormat = "mmm-yy"
$ws.Range("L29").HorizontalAlignment = -4108
$ws.Range("O29").Value = 1
$ws.Range("P29").Value = 5

# Row 30 : 12-Apr-2020
$ws.Range("A30").Value = 43933
$ws.Range("A30").NumberFormat = "d-mmm-yy"
$ws.Range("A30").HorizontalAlignment = -4108
$ws.Range("B30").Value = 6
$ws.Range("C30").Value = 766
$ws.Range("E30").Value = "Nairobi(4),Mombasa,Siaya"
$ws.Range("D30").Value = "UAE(2)"
$ws.Range("L30").Value = "25-59"
$ws.Range("G30").Value = "Community(4), Imported(2)"
$ws.Range("F30").Value = 197
$ws.Range("H30").Value = 1
$ws.Range("I30").Value = 1

# Row 31 : 13-Apr-2020
$ws.Range("A31").Value = 43934
$ws.Range("A31").NumberFormat = "d-mmm-yy"
$ws.Range("A31").HorizontalAlignment = -4108
$ws.Range("B31").Value = 11
$ws.Range("C31").Value = 674
$ws.Range("D31").Value = "UAE(4)"
$ws.Range("L31").Value = "1-42."
$ws.Range("L31").NumberFormat = "mmm-yy"
$ws.Range("L31").HorizontalAlignment = -4108
$ws.Range("E31").Value = "Mandera(4),Mombasa(3),Nairobi(2),Nakuru, Machakos"
$ws.Range("G31").Value = "Community(7),Imported(4)"
$ws.Range("F31").Value = 208
$ws.Range("H31").Value = 15
$ws.Range("I31").Value = 1
$ws.Range("O31").Value = 5
$ws.Range("P31").Value = 6

# Row 32 : 14-Apr-2020
$ws.Range("A32").Value = 43935
$ws.Range("A32").NumberFormat = "d-mmm-yy"
$ws.Range("A32").HorizontalAlignment = -4108
$ws.Range("B32").Value = 8
$ws.Range("C32").Value = 694
$ws.Range("G32").Value = "Community(2), Imported(6)"
$ws.Range("D32").Value = "UAE(2),UK,Pakistan,Zambia,Comoros"
$ws.Range("E32").Value = "Nairobi(6),Siaya,Nakuru"
$ws.Range("F32").Value = 216
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0

# Row 33 : 15-Apr-2020
$ws.Range("A33").Value = 43936
$ws.Range("A33").NumberFormat = "d-mmm-yy"
$ws.Range("A33").HorizontalAlignment = -4108
$ws.Range("B33").Value = 9
$ws.Range("C33").Value = 803
$ws.Range("D33").Value = "None"
$ws.Range("G33").Value = "Community(9)"
$ws.Range("E33").Value = "Nairobi(5), Mombasa(4)"
$ws.Range("L33").Value = "9-69."
$ws.Range("L33").NumberFormat = "mmm-yy"
$ws.Range("L33").HorizontalAlignment = -4108
$ws.Range("F33").Value = 225
$ws.Range("H33").Value = 12
$ws.Range("I33").Value = 1

# --- Update the view to match where the author left off editing ---
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
$ws.Range("C33").Select()
